$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "model_5_5_20"
$ws.Cells.Item(2, 2).Value = 0.4437020344020186
$ws.Cells.Item(2, 3).Value = 0.1639128337983419
$ws.Cells.Item(2, 4).Value = -0.5025018283923479
$ws.Cells.Item(2, 5).Value = -0.1477052686655687
$ws.Cells.Item(2, 6).Value = 0.6156579256057739
$ws.Cells.Item(2, 7).Value = 0.7016465067863464
$ws.Cells.Item(2, 8).Value = 2.942265033721924
$ws.Cells.Item(2, 9).Value = 1.75605034828186

$ws.Cells.Item(3, 1).Value = "model_5_5_19"
$ws.Cells.Item(3, 2).Value = 0.4560040278245617
$ws.Cells.Item(3, 3).Value = 0.1945075140017579
$ws.Cells.Item(3, 4).Value = -0.4568386020987696
$ws.Cells.Item(3, 5).Value = -0.111322504978457
$ws.Cells.Item(3, 6).Value = 0.6020432710647583
$ws.Cells.Item(3, 7).Value = 0.6759713292121887
$ws.Cells.Item(3, 8).Value = 2.852845191955566
$ws.Cells.Item(3, 9).Value = 1.700382828712463

$ws.Cells.Item(4, 1).Value = "model_5_5_18"
$ws.Cells.Item(4, 2).Value = 0.4625751248378341
$ws.Cells.Item(4, 3).Value = 0.2051591350178826
$ws.Cells.Item(4, 4).Value = -0.4233433843612238
$ws.Cells.Item(4, 5).Value = -0.08806084387972901
$ws.Cells.Item(4, 6).Value = 0.5947709679603577
$ws.Cells.Item(4, 7).Value = 0.6670324802398682
$ws.Cells.Item(4, 8).Value = 2.787253379821777
$ws.Cells.Item(4, 9).Value = 1.664791226387024

$ws.Cells.Item(5, 1).Value = "model_5_5_17"
$ws.Cells.Item(5, 2).Value = 0.4695370056520793
$ws.Cells.Item(5, 3).Value = 0.2139996192217355
$ws.Cells.Item(5, 4).Value = -0.3869481617171986
$ws.Cells.Item(5, 5).Value = -0.06357989706469969
$ws.Cells.Item(5, 6).Value = 0.5870662331581116
$ws.Cells.Item(5, 7).Value = 0.6596134901046753
$ws.Cells.Item(5, 8).Value = 2.715982675552368
$ws.Cells.Item(5, 9).Value = 1.627333998680115

$ws.Cells.Item(6, 1).Value = "model_5_5_16"
$ws.Cells.Item(6, 2).Value = 0.4890361506169921
$ws.Cells.Item(6, 3).Value = 0.2441349376097992
$ws.Cells.Item(6, 4).Value = -0.3106142436801993
$ws.Cells.Item(6, 5).Value = -0.008850358049641072
$ws.Cells.Item(6, 6).Value = 0.5654863715171814
$ws.Cells.Item(6, 7).Value = 0.6343238353729248
$ws.Cells.Item(6, 8).Value = 2.566502094268799
$ws.Cells.Item(6, 9).Value = 1.543594837188721

$ws.Cells.Item(7, 1).Value = "model_5_5_15"
$ws.Cells.Item(7, 2).Value = 0.4991155308119957
$ws.Cells.Item(7, 3).Value = 0.245337015652618
$ws.Cells.Item(7, 4).Value = -0.2570894845578111
$ws.Cells.Item(7, 5).Value = 0.02375066583830143
$ws.Cells.Item(7, 6).Value = 0.5543315410614014
$ws.Cells.Item(7, 7).Value = 0.6333150863647461
$ws.Cells.Item(7, 8).Value = 2.461687803268433
$ws.Cells.Item(7, 9).Value = 1.49371349811554

$ws.Cells.Item(8, 1).Value = "model_5_5_14"
$ws.Cells.Item(8, 2).Value = 0.5138247052748552
$ws.Cells.Item(8, 3).Value = 0.2498807592286936
$ws.Cells.Item(8, 4).Value = -0.186274280018131
$ws.Cells.Item(8, 5).Value = 0.06771001646738983
$ws.Cells.Item(8, 6).Value = 0.5380527377128601
$ws.Cells.Item(8, 7).Value = 0.62950199842453
$ws.Cells.Item(8, 8).Value = 2.323014259338379
$ws.Cells.Item(8, 9).Value = 1.426453351974487

$ws.Cells.Item(9, 1).Value = "model_5_5_13"
$ws.Cells.Item(9, 2).Value = 0.5308304599066355
$ws.Cells.Item(9, 3).Value = 0.2583412027022827
$ws.Cells.Item(9, 4).Value = -0.1070804175665383
$ws.Cells.Item(9, 5).Value = 0.1178563375921159
$ws.Cells.Item(9, 6).Value = 0.5192323923110962
$ws.Cells.Item(9, 7).Value = 0.6224019527435303
$ws.Cells.Item(9, 8).Value = 2.167933464050293
$ws.Cells.Item(9, 9).Value = 1.349726796150208

$ws.Cells.Item(10, 1).Value = "model_5_5_12"
$ws.Cells.Item(10, 2).Value = 0.5527359680699218
$ws.Cells.Item(10, 3).Value = 0.2642190410522356
$ws.Cells.Item(10, 4).Value = -0.007012162702188407
$ws.Cells.Item(10, 5).Value = 0.179845217086395
$ws.Cells.Item(10, 6).Value = 0.4949894845485687
$ws.Cells.Item(10, 7).Value = 0.6174692511558533
$ws.Cells.Item(10, 8).Value = 1.971975326538086
$ws.Cells.Item(10, 9).Value = 1.254880547523499

$ws.Cells.Item(11, 1).Value = "model_5_5_10"
$ws.Cells.Item(11, 2).Value = 0.5570710810396261
$ws.Cells.Item(11, 3).Value = 0.2284719890635835
$ws.Cells.Item(11, 4).Value = 0.06179633836513954
$ws.Cells.Item(11, 5).Value = 0.2109082438867399
$ws.Cells.Item(11, 6).Value = 0.490191787481308
$ws.Cells.Item(11, 7).Value = 0.6474682688713074
$ws.Cells.Item(11, 8).Value = 1.837231397628784
$ws.Cells.Item(11, 9).Value = 1.207352519035339

$ws.Cells.Item(12, 1).Value = "model_5_5_21"
$ws.Cells.Item(12, 2).Value = 0.5583351645460298
$ws.Cells.Item(12, 3).Value = -0.2056885621503861
$ws.Cells.Item(12, 4).Value = -0.06568409108182105
$ws.Cells.Item(12, 5).Value = 0.008058978027973618
$ws.Cells.Item(12, 6).Value = 0.4887928366661072
$ws.Cells.Item(12, 7).Value = 1.0118168592453
$ws.Cells.Item(12, 8).Value = 2.086869239807129
$ws.Cells.Item(12, 9).Value = 1.517722725868225

$ws.Cells.Item(13, 1).Value = "model_5_5_22"
$ws.Cells.Item(13, 2).Value = 0.5626389100006616
$ws.Cells.Item(13, 3).Value = -0.1501135590352825
$ws.Cells.Item(13, 4).Value = -0.0789160853437374
$ws.Cells.Item(13, 5).Value = 0.01622649083408434
$ws.Cells.Item(13, 6).Value = 0.4840298295021057
$ws.Cells.Item(13, 7).Value = 0.9651781916618347
$ws.Cells.Item(13, 8).Value = 2.112780570983887
$ws.Cells.Item(13, 9).Value = 1.505226016044617

$ws.Cells.Item(14, 1).Value = "model_5_5_9"
$ws.Cells.Item(14, 2).Value = 0.5634358415888661
$ws.Cells.Item(14, 3).Value = 0.2143277020570798
$ws.Cells.Item(14, 4).Value = 0.1124806628758737
$ws.Cells.Item(14, 5).Value = 0.23732339683173
$ws.Cells.Item(14, 6).Value = 0.4831478595733643
$ws.Cells.Item(14, 7).Value = 0.6593382358551025
$ws.Cells.Item(14, 8).Value = 1.737979292869568
$ws.Cells.Item(14, 9).Value = 1.166935920715332

$ws.Cells.Item(15, 1).Value = "model_5_5_11"
$ws.Cells.Item(15, 2).Value = 0.5634780784555854
$ws.Cells.Item(15, 3).Value = 0.2704565757809917
$ws.Cells.Item(15, 4).Value = 0.04933563697733168
$ws.Cells.Item(15, 5).Value = 0.2155954279668444
$ws.Cells.Item(15, 6).Value = 0.4831011593341827
$ws.Cells.Item(15, 7).Value = 0.6122347116470337
$ws.Cells.Item(15, 8).Value = 1.861632704734802
$ws.Cells.Item(15, 9).Value = 1.200180888175964

$ws.Cells.Item(16, 1).Value = "model_5_5_8"
$ws.Cells.Item(16, 2).Value = 0.56796453343292
$ws.Cells.Item(16, 3).Value = 0.2005386000255994
$ws.Cells.Item(16, 4).Value = 0.1567230377261105
$ws.Cells.Item(16, 5).Value = 0.2599684393121864
$ws.Cells.Item(16, 6).Value = 0.4781359136104584
$ws.Cells.Item(16, 7).Value = 0.6709100008010864
$ws.Cells.Item(16, 8).Value = 1.651341795921326
$ws.Cells.Item(16, 9).Value = 1.132287859916687

$ws.Cells.Item(17, 1).Value = "model_5_5_0"
$ws.Cells.Item(17, 2).Value = 0.5688575100333072
$ws.Cells.Item(17, 3).Value = 0.1018919747933875
$ws.Cells.Item(17, 4).Value = 0.4000017960718766
$ws.Cells.Item(17, 5).Value = 0.377847229516581
$ws.Cells.Item(17, 6).Value = 0.4771477282047272
$ws.Cells.Item(17, 7).Value = 0.7536945343017578
$ws.Cells.Item(17, 8).Value = 1.174942851066589
$ws.Cells.Item(17, 9).Value = 0.9519269466400146

$ws.Cells.Item(18, 1).Value = "model_5_5_7"
$ws.Cells.Item(18, 2).Value = 0.571715362281616
$ws.Cells.Item(18, 3).Value = 0.189754038856524
$ws.Cells.Item(18, 4).Value = 0.1975337060022163
$ws.Cells.Item(18, 5).Value = 0.2814139053082312
$ws.Cells.Item(18, 6).Value = 0.4739848971366882
$ws.Cells.Item(18, 7).Value = 0.6799604892730713
$ws.Cells.Item(18, 8).Value = 1.571424603462219
$ws.Cells.Item(18, 9).Value = 1.099475026130676

$ws.Cells.Item(19, 1).Value = "model_5_5_6"
$ws.Cells.Item(19, 2).Value = 0.572969649704502
$ws.Cells.Item(19, 3).Value = 0.1667621341157163
$ws.Cells.Item(19, 4).Value = 0.2349120486690038
$ws.Cells.Item(19, 5).Value = 0.297250601677016
$ws.Cells.Item(19, 6).Value = 0.4725967347621918
$ws.Cells.Item(19, 7).Value = 0.6992553472518921
$ws.Cells.Item(19, 8).Value = 1.498228669166565
$ws.Cells.Item(19, 9).Value = 1.075244069099426

$ws.Cells.Item(20, 1).Value = "model_5_5_4"
$ws.Cells.Item(20, 2).Value = 0.5730617473283708
$ws.Cells.Item(20, 3).Value = 0.1264596225066438
$ws.Cells.Item(20, 4).Value = 0.3023048497550352
$ws.Cells.Item(20, 5).Value = 0.3261364029475776
$ws.Cells.Item(20, 6).Value = 0.4724948108196259
$ws.Cells.Item(20, 7).Value = 0.7330772876739502
$ws.Cells.Item(20, 8).Value = 1.366257309913635
$ws.Cells.Item(20, 9).Value = 1.031047224998474

$ws.Cells.Item(21, 1).Value = "model_5_5_5"
$ws.Cells.Item(21, 2).Value = 0.5744110493216003
$ws.Cells.Item(21, 3).Value = 0.1481735283883039
$ws.Cells.Item(21, 4).Value = 0.2725844367299028
$ws.Cells.Item(21, 5).Value = 0.3145457994684896
$ws.Cells.Item(21, 6).Value = 0.4710015654563904
$ws.Cells.Item(21, 7).Value = 0.7148549556732178
$ws.Cells.Item(21, 8).Value = 1.42445695400238
$ws.Cells.Item(21, 9).Value = 1.048781514167786

$ws.Cells.Item(22, 1).Value = "model_5_5_23"
$ws.Cells.Item(22, 2).Value = 0.5783592143252647
$ws.Cells.Item(22, 3).Value = 0.07941049139846923
$ws.Cells.Item(22, 4).Value = -0.1316231137366113
$ws.Cells.Item(22, 5).Value = 0.05112734728463186
$ws.Cells.Item(22, 6).Value = 0.4666320979595184
$ws.Cells.Item(22, 7).Value = 0.7725610733032227
$ws.Cells.Item(22, 8).Value = 2.215993881225586
$ws.Cells.Item(22, 9).Value = 1.451825737953186

$ws.Cells.Item(23, 1).Value = "model_5_5_1"
$ws.Cells.Item(23, 2).Value = 0.5808566480439026
$ws.Cells.Item(23, 3).Value = 0.1402280814502024
$ws.Cells.Item(23, 4).Value = 0.4019968951279553
$ws.Cells.Item(23, 5).Value = 0.3901782009903616
$ws.Cells.Item(23, 6).Value = 0.4638682007789612
$ws.Cells.Item(23, 7).Value = 0.7215227484703064
$ws.Cells.Item(23, 8).Value = 1.171035885810852
$ws.Cells.Item(23, 9).Value = 0.9330598711967468

$ws.Cells.Item(24, 1).Value = "model_5_5_2"
$ws.Cells.Item(24, 2).Value = 0.5931181386068591
$ws.Cells.Item(24, 3).Value = 0.2000647665892772
$ws.Cells.Item(24, 4).Value = 0.3946328513689588
$ws.Cells.Item(24, 5).Value = 0.40312068519137
$ws.Cells.Item(24, 6).Value = 0.4502983093261719
$ws.Cells.Item(24, 7).Value = 0.6713076829910278
$ws.Cells.Item(24, 8).Value = 1.185456395149231
$ws.Cells.Item(24, 9).Value = 0.9132571816444397

$ws.Cells.Item(25, 1).Value = "model_5_5_3"
$ws.Cells.Item(25, 2).Value = 0.5955639851375476
$ws.Cells.Item(25, 3).Value = 0.2094891028499937
$ws.Cells.Item(25, 4).Value = 0.3733662940089891
$ws.Cells.Item(25, 5).Value = 0.3930466055519404
$ws.Cells.Item(25, 6).Value = 0.4475915133953094
$ws.Cells.Item(25, 7).Value = 0.6633987426757812
$ws.Cells.Item(25, 8).Value = 1.227101564407349
$ws.Cells.Item(25, 9).Value = 0.92867112159729

$ws.Cells.Item(26, 1).Value = "model_5_5_24"
$ws.Cells.Item(26, 2).Value = 0.8234903942607884
$ws.Cells.Item(26, 3).Value = 0.6900244314348998
$ws.Cells.Item(26, 4).Value = 0.4573846131001086
$ws.Cells.Item(26, 5).Value = 0.5831836203725356
$ws.Cells.Item(26, 6).Value = 0.1953441202640533
$ws.Cells.Item(26, 7).Value = 0.2601322829723358
$ws.Cells.Item(26, 8).Value = 1.062573194503784
$ws.Cells.Item(26, 9).Value = 0.6377512812614441
